# Generate Report for Handoff
# Adds two new handed-off files (a24594a0-...md and c141966c-...md) as new
# rows across the Overview, zh-cn and de-de sheets, mirroring the existing
# rows for 0bf4a485-...md / 17ad62a6-...md.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet (A1:G3 -> A1:G5)
#   A File Name | B Path And Name | C Extension | D Publish URL
#   E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------

$wsOverview.Range("A4").Value = "a24594a0-b84d-412c-8b49-21697e41ef2a.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-28 20:40:41"

$wsOverview.Range("A5").Value = "c141966c-2899-4b3e-9e77-cedc06e759c9.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-28 20:40:41"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a24594a0b84d412c8b4921697e41ef2a0000000/e2e/a24594a0-b84d-412c-8b49-21697e41ef2a.md",
    "", "",
    "e2e\a24594a0-b84d-412c-8b49-21697e41ef2a.md"
)
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B5"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c141966c28994b3e9e77cedc06e759c900000000/e2e/c141966c-2899-4b3e-9e77-cedc06e759c9.md",
    "", "",
    "e2e\c141966c-2899-4b3e-9e77-cedc06e759c9.md"
)

# ---------------------------------------------------------------------------
# zh-cn sheet (A1:P3 -> A1:P5)
#   A Source File Name | B File Extension | C Status | D Source Path
#   E Priority | F Content Duplicate | G Latest Handoff File
#   H Latest Handoff Datetime | I Latest Target File | J Latest Handback File
#   K Latest Handback DateTime | L Reference Tokens | M To be localized
#   N Dependency From | O Has metadata | P Error Detail
# ---------------------------------------------------------------------------

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "a24594a0-b84d-412c-8b49-21697e41ef2a.11fd45dfa5f21bf61530f55eb23ac5c65ab99937.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-28 20:40:35"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "c141966c-2899-4b3e-9e77-cedc06e759c9.587f99387ca59b7b5f80d317ee13b03c3866ef6f.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-28 20:40:35"
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = ""
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L5").Value = ""
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("N5").Value = ""
$wsZhCn.Range("O5").Value = "False"
$wsZhCn.Range("P5").Value = ""

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a24594a0b84d412c8b4921697e41ef2a0000000/e2e/a24594a0-b84d-412c-8b49-21697e41ef2a.md",
    "", "",
    "a24594a0-b84d-412c-8b49-21697e41ef2a.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A5"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c141966c28994b3e9e77cedc06e759c900000000/e2e/c141966c-2899-4b3e-9e77-cedc06e759c9.md",
    "", "",
    "c141966c-2899-4b3e-9e77-cedc06e759c9.md"
)

# ---------------------------------------------------------------------------
# de-de sheet (A1:P3 -> A1:P5) - same columns as zh-cn
# ---------------------------------------------------------------------------

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "a24594a0-b84d-412c-8b49-21697e41ef2a.11fd45dfa5f21bf61530f55eb23ac5c65ab99937.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-28 20:40:41"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "c141966c-2899-4b3e-9e77-cedc06e759c9.587f99387ca59b7b5f80d317ee13b03c3866ef6f.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-28 20:40:41"
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = ""
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L5").Value = ""
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("N5").Value = ""
$wsDeDe.Range("O5").Value = "False"
$wsDeDe.Range("P5").Value = ""

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a24594a0b84d412c8b4921697e41ef2a0000000/e2e/a24594a0-b84d-412c-8b49-21697e41ef2a.md",
    "", "",
    "a24594a0-b84d-412c-8b49-21697e41ef2a.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A5"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c141966c28994b3e9e77cedc06e759c900000000/e2e/c141966c-2899-4b3e-9e77-cedc06e759c9.md",
    "", "",
    "c141966c-2899-4b3e-9e77-cedc06e759c9.md"
)

# ---------------------------------------------------------------------------
# Table ranges grow from 3 data rows to 5; resize each ListObject to match
# the new dimensions (A1:G5 for Overview, A1:P5 for zh-cn / de-de).
# ---------------------------------------------------------------------------

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G5"))
$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P5"))
$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P5"))

Write-Host "Handback rows added for a24594a0-b84d-412c-8b49-21697e41ef2a.md and c141966c-2899-4b3e-9e77-cedc06e759c9.md"
